# Keys_Onboarding / TestData.xlsx — add a new data row (row 3) to the
# LoginPage sheet with a second login/registration record, including its
# two hyperlinks (URL + mailto), mirroring the existing row 2 pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new data row -----------------------------------------------------
$ws.Range("A3").Value = "http://new-keys.azurewebsites.net/Account/Login"
$ws.Range("B3").Value = "baha.godbole@gmail.com"
$ws.Range("C3").Value = "abc1abc1"
$ws.Range("D3").Value = "description added"

# --- hyperlinks (B3 first, then A3 — matches author order) -----------
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:baha.godbole@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A3"), "http://new-keys.azurewebsites.net/Account/Login")

# --- re-apply row 2's Hyperlink formatting so A3/B3 share its style --
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)

# --- leave the selection where the author's cursor ended up ----------
[void]$ws.Range("D4").Select()
